# Update quantities / amounts on the "Bill Summary" sheet.
# Plain numeric cells (column C) are written directly as numbers.
# Text-formatted numeric-looking amounts (column G / H, stored as text in the
# original workbook) are written with a leading apostrophe so Excel keeps
# them as text instead of re-interpreting them as numbers, then
# ClearFormats() strips the "quote prefix" style Excel tags onto the cell so
# the cell keeps using the workbook's default (unstyled) format exactly like
# every other untouched cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $text) {
    $rng = $ws.Range($rangeAddress)
    $rng.Value = "'" + $text
    $rng.ClearFormats()
}

# Row 8
$ws.Range("C8").Value = 66

# Row 9
$ws.Range("C9").Value = 30
Set-TextValue "G9" "7680.00"

# Row 10
$ws.Range("C10").Value = 20
Set-TextValue "G10" "9440.00"

# Row 11
$ws.Range("C11").Value = 8
Set-TextValue "G11" "5296.00"

# Row 12
$ws.Range("C12").Value = 98

# Row 13
$ws.Range("C13").Value = 76
Set-TextValue "G13" "10336.00"

# Row 14
$ws.Range("C14").Value = 11
Set-TextValue "G14" "253.00"

# Row 15
$ws.Range("C15").Value = 17

# Row 16
$ws.Range("C16").Value = 8

# Row 17
$ws.Range("C17").Value = 81

# Row 19 - Grand Total Rs.
Set-TextValue "G19" "33005.00"
Set-TextValue "H19" "33005.00"

# Row 21 - NET PAYABLE AMOUNT Rs.
Set-TextValue "G21" "33005.00"
Set-TextValue "H21" "33005.00"
